$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "Leo Parisi "
$ws.Range("B46").Value = "Daniele Dalbosco | IMONTAGNA"
$ws.Range("C46").Value = "Leonardo Viola | SHARK ATTACK"
$ws.Range("D46").Value = "Geremia  Carollo | FC SAVIGNANO"
$ws.Range("E46").Value = "Luca Frasca | Clitoriders"
$ws.Range("F46").Value = "Moris Benedetti | Gli Introvabili"
